$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "318.66"
Set-TextValue "E2" "3.96%"
Set-TextValue "D3" "39.72"
Set-TextValue "E3" "1.68%"
Set-TextValue "D4" "5.139"
Set-TextValue "E4" "0.83%"
Set-TextValue "D5" "0.08222"
Set-TextValue "E5" "2.14%"
Set-TextValue "D6" "2.068"
Set-TextValue "E6" "5.67%"
Set-TextValue "D7" "8.312"
Set-TextValue "E7" "3.83%"
Set-TextValue "B8" "GateToken"
Set-TextValue "C8" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D8" "4.297"
Set-TextValue "E8" "2.43%"
Set-TextValue "B9" "MXToken"
Set-TextValue "C9" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D9" "0.9348"
Set-TextValue "E9" "0.31%"
Set-TextValue "B10" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C10" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D10" "0.1357"
Set-TextValue "E10" "-6.16%"
Set-TextValue "B11" "WazirX"
Set-TextValue "C11" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1984"
Set-TextValue "E11" "2.72%"
Set-TextValue "B12" "MandalaExchangeToken"
Set-TextValue "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.09074"
Set-TextValue "E12" "-0.40%"
Set-TextValue "B13" "BitrueCoin"
Set-TextValue "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.03505"
Set-TextValue "E13" "-0.05%"
Set-TextValue "B14" "BitMartToken"
Set-TextValue "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09830"
Set-TextValue "E14" "0.41%"
Set-TextValue "B15" "BitForexToken"
Set-TextValue "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D15" "0.001398"
Set-TextValue "E15" "-0.15%"
Set-TextValue "B16" "TigerCash"
Set-TextValue "C16" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D16" "0.006290"
Set-TextValue "E16" "4.20%"
Set-TextValue "B17" "LEO"
Set-TextValue "C17" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D17" "3.683"
Set-TextValue "E17" "-2.47%"
Set-TextValue "E18" "-2.66%"
Set-TextValue "D19" "0.3474"
Set-TextValue "E19" "1.49%"
Set-TextValue "D20" "0.1294"
Set-TextValue "E20" "-3.93%"
Set-TextValue "D21" "4.901"
Set-TextValue "E21" "7.28%"
Set-TextValue "D22" "0.2447"
Set-TextValue "E22" "1.25%"
Set-TextValue "D23" "0.04319"
Set-TextValue "E23" "-1.50%"
Set-TextValue "D24" "0.001226"
Set-TextValue "E24" "-0.87%"
Set-TextValue "E25" "11.66%"
Set-TextValue "D26" "0.0001298"
Set-TextValue "E26" "-0.30%"
Set-TextValue "D27" "0.0003993"
Set-TextValue "E27" "-10.21%"
Set-TextValue "D39" "0.02217"
Set-TextValue "E39" "9.03%"
Set-TextValue "D40" "0.05228"
Set-TextValue "E40" "3.47%"
Set-TextValue "D41" "0.007650"
Set-TextValue "E41" "2.76%"
Set-TextValue "D42" "0.009738"
Set-TextValue "E42" "-5.41%"
Set-TextValue "D43" "0.1389"
Set-TextValue "E43" "3.30%"
Set-TextValue "D44" "0.002091"
Set-TextValue "E44" "-1.54%"
Set-TextValue "D45" "0.009195"
Set-TextValue "E45" "0.85%"
Set-TextValue "D46" "0.00006548"
Set-TextValue "E46" "5.53%"
Set-TextValue "D47" "0.00000000749"
Set-TextValue "E47" "-0.34%"
Set-TextValue "D48" "0.002985"
Set-TextValue "E48" "-3.63%"
Set-TextValue "D49" "0.001687"
Set-TextValue "E49" "5.35%"
Set-TextValue "D50" "0.00002097"
Set-TextValue "E50" "-0.34%"
Set-TextValue "D51" "0.0001997"
Set-TextValue "E51" "-0.34%"
